$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.883.77'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '3.280.27'
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''584.92'
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").Value = '''180.10'
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("D7").Value = '''0.653'
$ws.Range("E7").Value = '  +8.66%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -3.41%  '

$ws.Range("E10").Value = '  +2.07%  '

$ws.Range("D11").Value = '''0.405'
$ws.Range("E11").Value = '  +0.42%  '

$ws.Range("D12").Value = '3.848.92'
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("E13").Value = '  -4.53%  '

$ws.Range("D14").Value = '65.963.95'
$ws.Range("E14").Value = '  -1.09%  '

$ws.Range("D15").Value = '''26.33'
$ws.Range("E15").Value = '  -3.27%  '

$ws.Range("E16").Value = '  -2.44%  '

$ws.Range("D17").Value = '3.233.67'
$ws.Range("E17").Value = '  -2.33%  '

$ws.Range("D18").Value = '''429.61'
$ws.Range("E18").Value = '  -1.66%  '

$ws.Range("D19").Value = '''13.24'
$ws.Range("E19").Value = '  -3.82%  '

$ws.Range("D20").Value = '''5.52'
$ws.Range("E20").Value = '  -2.95%  '

$ws.Range("E21").Value = '  -3.30%  '

$ws.Range("D22").Value = '''72.19'
$ws.Range("E22").Value = '  -2.27%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '''5.68'
$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").Value = '3.433.71'
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("E26").Value = '  -1.17%  '

$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '''0.0000113'
$ws.Range("E27").Value = '  -4.44%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.196'
$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("D29").Value = '''8.91'
$ws.Range("E29").Value = '  -1.82%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").Value = '''1.97'
$ws.Range("E31").Value = '  +0.49%  '

$ws.Range("D32").Value = '''22.31'
$ws.Range("E32").Value = '  -2.57%  '

$ws.Range("D33").Value = '''1.00'

$ws.Range("D34").Value = '''5.17'
$ws.Range("E34").Value = '  -2.98%  '

$ws.Range("D35").Value = '''6.60'
$ws.Range("E35").Value = '  -2.68%  '

$ws.Range("E36").Value = '  -3.35%  '

$ws.Range("D37").Value = '''158.78'
$ws.Range("E37").Value = '  -0.89%  '

$ws.Range("D38").Value = '''1.42'
$ws.Range("E38").Value = '  -5.70%  '

$ws.Range("D39").Value = '''26.48'
$ws.Range("E39").Value = '  -3.49%  '

$ws.Range("E40").Value = '  -4.11%  '

$ws.Range("D41").Value = '2.780.80'
$ws.Range("E41").Value = '  -0.99%  '

$ws.Range("D42").Value = '''0.768'
$ws.Range("E42").Value = '  -2.96%  '

$ws.Range("D43").Value = '''4.33'
$ws.Range("E43").Value = '  -3.14%  '

$ws.Range("D44").Value = '''40.06'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").Value = '''0.0659'
$ws.Range("E45").Value = '  -2.75%  '

$ws.Range("D46").Value = '''5.91'
$ws.Range("E46").Value = '  -5.40%  '

$ws.Range("E47").Value = '  -1.84%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '''314.94'
$ws.Range("E48").Value = '  -1.20%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''23.16'
$ws.Range("E49").Value = '  -4.76%  '

$ws.Range("E50").Value = '  -2.09%  '

$ws.Range("E51").Value = '  +6.43%  '
